$p = $ppt.ActivePresentation

# --------------------------------------------------------------------------
# Slide 1 - intro paragraph textbox ("TextBox 3")
# --------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$tb1 = $s1.Shapes.Item(3).TextFrame.TextRange

# Replace the whole first sentence-run: "they system will always be
# appreciated" -> "the system are appreciated" (keep identical run bounds
# so the run does not get split).
$run1 = $tb1.Characters(1, 282)
$run1.Text = "The quick response homework system is being developed to give you relevant practice on numerical concepts.  This system is free to use and is still under construction so suggestions on improving the system are appreciated.  Unlike systems you may have used in the past, "

# Merge the trailing three runs ("and work it...", "a textbook.  ", "You can
# check you answer...") into a single run, fixing "you answer" -> "your
# answers" at the same time.
$fullNow = $tb1.Text
$idx = $fullNow.IndexOf("and work it as you would a problem from ")
$tailLen = $tb1.Length - $idx
$tailRun = $tb1.Characters($idx + 1, $tailLen)
$tailRun.Text = "and work it as you would a problem from a textbook.  You can check your answers on your phone or computer.  If you get stuck, you may try working the base-case problem.  You are provided the answers to the basecase problem.  Questions also come with reflections that your instructor may assign.  "

# --------------------------------------------------------------------------
# Slide 2 - "Please Note this" callout (single run)
# --------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$noteRange = $s2.Shapes.Item(5).TextFrame.TextRange
$noteRange.Text = "Please Note this – if it the problem does not fully load it will have markups like ##mdot,num,20## in the statement"

# Slide 2 - "In your Browser..." heading (only the first run changes)
$browserRange = $s2.Shapes.Item(6).TextFrame.TextRange
$browserFirst = $browserRange.Characters(1, 46)
$browserFirst.Text = "In your Browser (Chrome is recommended) Type:  "

# --------------------------------------------------------------------------
# Slide 3 - "Ctrl P" title: split the trailing run into 3 runs so "Ctrl P"
# gets its own red/underlined formatting.
# --------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$titleRange = $s3.Shapes.Item(2).TextFrame.TextRange

$part2 = $titleRange.Characters(91, 6)
$part2.Text = "Ctrl P"
$part2.Font.Underline = $true
$part2.Font.Color.RGB = 255

# --------------------------------------------------------------------------
# Slide 4 - "The Base Case..." title -> "The Base-Case..."
# --------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$baseCaseRange = $s4.Shapes.Item(1).TextFrame.TextRange
$baseCaseFirst = $baseCaseRange.Characters(1, 42)
$baseCaseFirst.Text = "The Base-Case is the Same for All Students"

# --------------------------------------------------------------------------
# Slide 6 - "You can check each part..." callout (single run)
# --------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$checkRange = $s6.Shapes.Item(3).TextFrame.TextRange
$checkRange.Text = "You can check each part of your answer as many times as you like without penalty – however the system will slow after several incorrect tries"

# --------------------------------------------------------------------------
# Slide 6 - duplicate the red left-arrow, unrotated, pointing at the new
# "Get answer to Base-Case Problem Here" textbox lower on the slide.
# --------------------------------------------------------------------------
$origArrow = $s6.Shapes.Item(5)
$newArrow = $origArrow.Duplicate()
$newArrow = $s6.Shapes.Item($s6.Shapes.Count)
$newArrow.Name = "Arrow: Left 6"
$newArrow.Rotation = 0
$newArrow.Left = 363.6225280761719
$newArrow.Top = 441.4634094238281
